# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
#
# Column D ("Price") and E ("Volume(1h)") are stored as plain text in this
# sheet (e.g. "603.54", "  +0.10%  ") rather than numbers, so that values
# like "1.00" or "0.999" keep their exact printed form. Assigning a
# number-looking string straight to Range.Value makes Excel auto-convert the
# cell to a real number (same as typing it in manually), which would both
# change the cell's type and normalise away formatting like trailing zeros.
# Set-TextValue avoids that: it stamps the cell as Text ("@") just long
# enough to accept the literal string, then restores the cell's original
# style so no formatting/style footprint is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $ws.Range($cell).Style
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($cell).Style = $origStyle
}

$ws.Range("D2").Value = "66.492.09"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "3.183.40"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "603.50"
$ws.Range("E5").Value = "  +0.11%  "
Set-TextValue "D6" "155.83"
$ws.Range("E6").Value = "  +2.27%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +2.96%  "
$ws.Range("D9").Value = "3.181.54"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("E10").Value = "  -1.82%  "
Set-TextValue "D11" "5.70"
$ws.Range("E11").Value = "  -7.54%  "
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("E13").Value = "  -1.47%  "
Set-TextValue "D14" "38.93"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "3.705.65"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "66.522.54"
$ws.Range("E16").Value = "  +0.59%  "
Set-TextValue "D17" "7.46"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "3.181.74"
$ws.Range("E18").Value = "  -1.25%  "
Set-TextValue "D20" "513.58"
$ws.Range("E20").Value = "  +0.10%  "
Set-TextValue "D21" "15.56"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("E22").Value = "  -0.32%  "
Set-TextValue "D23" "8.24"
$ws.Range("E23").Value = "  +3.01%  "
Set-TextValue "D24" "14.97"
$ws.Range("E24").Value = "  -1.57%  "
Set-TextValue "D25" "84.80"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.72%  "
Set-TextValue "D28" "9.19"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  +7.32%  "
Set-TextValue "D30" "3.09"
$ws.Range("E30").Value = "  +6.89%  "
Set-TextValue "D31" "7.00"
$ws.Range("E31").Value = "  +2.05%  "
Set-TextValue "D32" "28.19"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -1.92%  "
Set-TextValue "D35" "6.57"
$ws.Range("E35").Value = "  -1.07%  "
Set-TextValue "D36" "514.74"
$ws.Range("E36").Value = "  +5.59%  "
Set-TextValue "D37" "54.91"
$ws.Range("E37").Value = "  -1.35%  "
Set-TextValue "D38" "0.0896"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("E39").Value = "  +0.07%  "
$ws.Range("E40").Value = "  +5.55%  "
Set-TextValue "D41" "8.93"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "0.0₃0685"
$ws.Range("E42").Value = "  +5.97%  "
Set-TextValue "D43" "0.301"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("E44").Value = "  -6.79%  "
Set-TextValue "D45" "2.44"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").Value = "2.856.78"
$ws.Range("E46").Value = "  -5.67%  "
Set-TextValue "D47" "28.37"
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("E48").Value = "  +2.96%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.117"
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D50" "0.999"
$ws.Range("E50").Value = "  -0.12%  "
Set-TextValue "D51" "2.62"
$ws.Range("E51").Value = "  +7.40%  "
